# "Fruta / hortaliza, semanal" — weekly refresh of the Sandia / Vega Modelo de
# Temuco subconjunto: a new weekly snapshot (4 quality rows, date 44610) is
# inserted at the top of the data block, pushing the previously-existing
# 458-491 rows down by four (to 462-495).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the existing data block (old rows
# 458-491 shift down to 462-495, inheriting the row-458 formatting).
$ws.Range("A458:R461").EntireRow.Insert()

# Columns shared by every row in this subconjunto (market/category block).
$Mercado   = 10
$MercadoNm = "Vega Modelo de Temuco"
$Region    = "La Araucanía"
$Codreg    = 9
$CatId     = 100112028
$Categoria = "Sandia"
$Variedad  = "Sin especificar"
$KgUnid    = 1
$Clasif    = "Hortaliza"

# New weekly rows (one per calidad), date serial 44610.
$newRows = @(
    @{ Row = 458; Fecha = 44610; Calidad = "Extra";    Volumen = 2000; PMin = 2800; PMax = 3000; PProm = 2900; Unidad = "`$/unidad"; Origen = "Región del Maule"; PKg = 2900 },
    @{ Row = 459; Fecha = 44610; Calidad = "Primera";  Volumen = 4000; PMin = 2500; PMax = 2500; PProm = 2500; Unidad = "`$/unidad"; Origen = "Región del Maule"; PKg = 2500 },
    @{ Row = 460; Fecha = 44610; Calidad = "Segunda";  Volumen = 3000; PMin = 2000; PMax = 2000; PProm = 2000; Unidad = "`$/unidad"; Origen = "Región del Maule"; PKg = 2000 },
    @{ Row = 461; Fecha = 44610; Calidad = "Tercera";  Volumen = 1500; PMin = 1000; PMax = 1000; PProm = 1000; Unidad = "`$/unidad"; Origen = "Región del Maule"; PKg = 1000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value2 = $Mercado
    $ws.Range("B$row").Value2 = $MercadoNm
    $ws.Range("C$row").Value2 = $Region
    $ws.Range("D$row").Value2 = $r.Fecha
    $ws.Range("E$row").Value2 = $Codreg
    $ws.Range("F$row").Value2 = $CatId
    $ws.Range("G$row").Value2 = $Categoria
    $ws.Range("H$row").Value2 = $Variedad
    $ws.Range("I$row").Value2 = $r.Calidad
    $ws.Range("J$row").Value2 = $r.Volumen
    $ws.Range("K$row").Value2 = $r.PMin
    $ws.Range("L$row").Value2 = $r.PMax
    $ws.Range("M$row").Value2 = $r.PProm
    $ws.Range("N$row").Value2 = $r.Unidad
    $ws.Range("O$row").Value2 = $r.Origen
    $ws.Range("P$row").Value2 = $r.PKg
    $ws.Range("Q$row").Value2 = $KgUnid
    $ws.Range("R$row").Value2 = $Clasif
}
